$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2491.2817
$ws.Range("I15").Value = 2491.2817
$ws.Range("K15").Value = 7473.8451
$ws.Range("M15").Value = -7304.8451
$ws.Range("H70").Value = 1133.1912
$ws.Range("I70").Value = 962.58826
$ws.Range("J70").Value = 1303.7941
$ws.Range("K70").Value = 2887.76478
$ws.Range("L70").Value = 3911.3823
$ws.Range("M70").Value = -2617.76478
$ws.Range("N70").Value = -4451.3823
$ws.Range("H73").Value = 1133.1912
$ws.Range("I73").Value = 962.58826
$ws.Range("J73").Value = 1303.7941
$ws.Range("K73").Value = 2887.76478
$ws.Range("L73").Value = 3911.3823
$ws.Range("M73").Value = -1951.76478
$ws.Range("N73").Value = -5783.3823
$ws.Range("H76").Value = 3128.5715
$ws.Range("I76").Value = 2844.4443
$ws.Range("J76").Value = 3640
$ws.Range("K76").Value = 2844.4443
$ws.Range("L76").Value = 3640
$ws.Range("M76").Value = -2529.4443
$ws.Range("N76").Value = -4270
$ws.Range("H79").Value = 3128.5715
$ws.Range("I79").Value = 2844.4443
$ws.Range("J79").Value = 3640
$ws.Range("K79").Value = 2844.4443
$ws.Range("L79").Value = 3640
$ws.Range("M79").Value = -1752.4443
$ws.Range("N79").Value = -5824
$ws.Range("H86").Value = 66542.42
$ws.Range("I86").Value = 154925.75
$ws.Range("J86").Value = 2263.6365
$ws.Range("K86").Value = 154925.75
$ws.Range("L86").Value = 2263.6365
$ws.Range("M86").Value = -153802.75
$ws.Range("N86").Value = -4509.636500000001
$ws.Range("H89").Value = 66542.42
$ws.Range("I89").Value = 154925.75
$ws.Range("J89").Value = 2263.6365
$ws.Range("K89").Value = 774628.75
$ws.Range("L89").Value = 11318.1825
$ws.Range("M89").Value = -769012.75
$ws.Range("N89").Value = -22550.1825
$ws.Range("H98").Value = 2443.5715
$ws.Range("I98").Value = 1982.7273
$ws.Range("J98").Value = 4133.3335
$ws.Range("K98").Value = 1982.7273
$ws.Range("L98").Value = 4133.3335
$ws.Range("M98").Value = -484.7273
$ws.Range("N98").Value = -7129.3335
$ws.Range("H122").Value = 2443.5715
$ws.Range("I122").Value = 1982.7273
$ws.Range("J122").Value = 4133.3335
$ws.Range("K122").Value = 5948.1819
$ws.Range("L122").Value = 12400.0005
$ws.Range("M122").Value = -3498.1819
$ws.Range("N122").Value = -17300.0005
$ws.Range("H137").Value = 10418559
$ws.Range("I137").Value = 16668685
$ws.Range("J137").Value = 1682.5
$ws.Range("K137").Value = 50006055
$ws.Range("L137").Value = 5047.5
$ws.Range("M137").Value = -50003505
$ws.Range("N137").Value = -10147.5
$ws.Range("H138").Value = 2303.8867
$ws.Range("I138").Value = 2160.75
$ws.Range("J138").Value = 2390.6365
$ws.Range("K138").Value = 6482.25
$ws.Range("L138").Value = 7171.9095
$ws.Range("M138").Value = -1342.25
$ws.Range("N138").Value = -17451.9095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 52.25
$ws.Range("I5").Value = 52.25
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 52.25
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 59.75
$ws.Range("N5").ClearContents()
$ws.Range("H8").Value = 1970
$ws.Range("I8").Value = 980
$ws.Range("J8").Value = 2300
$ws.Range("K8").Value = 980
$ws.Range("L8").Value = 2300
$ws.Range("M8").Value = -836
$ws.Range("N8").Value = -2588
$ws.Range("H32").Value = 5062903.5
$ws.Range("I32").Value = 7017.8037
$ws.Range("J32").Value = 33375864
$ws.Range("K32").Value = 7017.8037
$ws.Range("L32").Value = 33375864
$ws.Range("M32").Value = -6730.8037
$ws.Range("N32").Value = -33376438
$ws.Range("H113").Value = 20398
$ws.Range("J113").Value = 20398
$ws.Range("L113").Value = 20398
$ws.Range("N113").Value = -29076

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 52.25
$ws.Range("I4").Value = 52.25
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 52.25
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 62.75
$ws.Range("N4").ClearContents()
$ws.Range("H86").Value = 1703.4615
$ws.Range("I86").Value = 1560.4546
$ws.Range("J86").Value = 2490
$ws.Range("K86").Value = 1560.4546
$ws.Range("L86").Value = 2490
$ws.Range("M86").Value = -437.4546
$ws.Range("N86").Value = -4736
$ws.Range("H89").Value = 1703.4615
$ws.Range("I89").Value = 1560.4546
$ws.Range("J89").Value = 2490
$ws.Range("K89").Value = 7802.273
$ws.Range("L89").Value = 12450
$ws.Range("M89").Value = -2186.273
$ws.Range("N89").Value = -23682
$ws.Range("H94").Value = 1125.25
$ws.Range("I94").Value = 1125.25
$ws.Range("K94").Value = 1125.25
$ws.Range("M94").Value = -674.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 6704.8667
$ws.Range("I7").Value = 9121.362999999999
$ws.Range("J7").Value = 59.5
$ws.Range("K7").Value = 9121.362999999999
$ws.Range("L7").Value = 59.5
$ws.Range("M7").Value = -9008.362999999999
$ws.Range("N7").Value = -285.5
$ws.Range("H31").Value = 4903937
$ws.Range("I31").Value = 1500.6389
$ws.Range("J31").Value = 10419178
$ws.Range("K31").Value = 1500.6389
$ws.Range("L31").Value = 10419178
$ws.Range("M31").Value = -1205.6389
$ws.Range("N31").Value = -10419768
$ws.Range("H34").Value = 4903937
$ws.Range("I34").Value = 1500.6389
$ws.Range("J34").Value = 10419178
$ws.Range("K34").Value = 1500.6389
$ws.Range("L34").Value = 10419178
$ws.Range("M34").Value = -1298.6389
$ws.Range("N34").Value = -10419582

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 683
$ws.Range("I5").Value = 427.13635
$ws.Range("J5").Value = 1245.9
$ws.Range("K5").Value = 1281.40905
$ws.Range("L5").Value = 3737.7
$ws.Range("M5").Value = -1169.40905
$ws.Range("N5").Value = -3961.7
$ws.Range("H134").Value = 3127.641
$ws.Range("I134").Value = 1181.9286
$ws.Range("J134").Value = 4217.24
$ws.Range("K134").Value = 3545.7858
$ws.Range("L134").Value = 12651.72
$ws.Range("M134").Value = 1524.2142
$ws.Range("N134").Value = -22791.72
$ws.Range("H135").Value = 683
$ws.Range("I135").Value = 427.13635
$ws.Range("J135").Value = 1245.9
$ws.Range("K135").Value = 3844.22715
$ws.Range("L135").Value = 11213.1
$ws.Range("M135").Value = -1309.22715
$ws.Range("N135").Value = -16283.1
$ws.Range("H136").Value = 4125.737
$ws.Range("I136").Value = 2998.889
$ws.Range("J136").Value = 5139.9
$ws.Range("K136").Value = 8996.667000000001
$ws.Range("L136").Value = 15419.7
$ws.Range("M136").Value = -3896.667000000001
$ws.Range("N136").Value = -25619.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2455.3572
$ws.Range("I80").Value = 2534.0908
$ws.Range("K80").Value = 2534.0908
$ws.Range("M80").Value = -1536.0908
$ws.Range("H83").Value = 2455.3572
$ws.Range("I83").Value = 2534.0908
$ws.Range("K83").Value = 12670.454
$ws.Range("M83").Value = -7678.454
$ws.Range("H132").Value = 2253.889
$ws.Range("I132").Value = 2814.7778
$ws.Range("J132").Value = 1973.4445
$ws.Range("K132").Value = 8444.3334
$ws.Range("L132").Value = 5920.333500000001
$ws.Range("M132").Value = -5914.3334
$ws.Range("N132").Value = -10980.3335
$ws.Range("H137").Value = 40000
$ws.Range("J137").Value = 40000
$ws.Range("L137").Value = 40000
$ws.Range("N137").Value = -50200

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3368.4614
$ws.Range("I68").Value = 3060
$ws.Range("J68").Value = 3632.8572
$ws.Range("K68").Value = 3060
$ws.Range("L68").Value = 3632.8572
$ws.Range("M68").Value = -2311
$ws.Range("N68").Value = -5130.8572
$ws.Range("H71").Value = 3368.4614
$ws.Range("I71").Value = 3060
$ws.Range("J71").Value = 3632.8572
$ws.Range("K71").Value = 15300
$ws.Range("L71").Value = 18164.286
$ws.Range("M71").Value = -11556
$ws.Range("N71").Value = -25652.286
$ws.Range("H82").Value = 3070
$ws.Range("I82").Value = 2433.3333
$ws.Range("J82").Value = 4980
$ws.Range("K82").Value = 2433.3333
$ws.Range("L82").Value = 4980
$ws.Range("M82").Value = -2072.3333
$ws.Range("N82").Value = -5702
$ws.Range("H85").Value = 3070
$ws.Range("I85").Value = 2433.3333
$ws.Range("J85").Value = 4980
$ws.Range("K85").Value = 2433.3333
$ws.Range("L85").Value = 4980
$ws.Range("M85").Value = -1185.3333
$ws.Range("N85").Value = -7476
$ws.Range("H93").Value = 186533.62
$ws.Range("I93").Value = 264360.75
$ws.Range("J93").Value = 1694.25
$ws.Range("K93").Value = 264360.75
$ws.Range("L93").Value = 1694.25
$ws.Range("M93").Value = -263112.75
$ws.Range("N93").Value = -4190.25
$ws.Range("H134").Value = 47400
$ws.Range("J134").Value = 47400
$ws.Range("L134").Value = 47400
$ws.Range("N134").Value = -57540
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 19231886
$ws.Range("I136").Value = 25000804
$ws.Range("J136").Value = 2161.5
$ws.Range("K136").Value = 75002412
$ws.Range("L136").Value = 6484.5
$ws.Range("M136").Value = -74999862
$ws.Range("N136").Value = -11584.5
$ws.Range("H137").Value = 31000
$ws.Range("J137").Value = 31000
$ws.Range("L137").Value = 31000
$ws.Range("N137").Value = -41200
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 7994.1665
$ws.Range("I122").Value = 9270.588
$ws.Range("J122").Value = 4894.2856
$ws.Range("K122").Value = 27811.764
$ws.Range("L122").Value = 14682.8568
$ws.Range("M122").Value = -25361.764
$ws.Range("N122").Value = -19582.8568
$ws.Range("H126").Value = 7199.875
$ws.Range("I126").Value = 7467.8184
$ws.Range("J126").Value = 4252.5
$ws.Range("K126").Value = 22403.4552
$ws.Range("L126").Value = 12757.5
$ws.Range("M126").Value = -19933.4552
$ws.Range("N126").Value = -17697.5
$ws.Range("H132").Value = 2639.9119
$ws.Range("I132").Value = 2901.818
$ws.Range("J132").Value = 2514.652
$ws.Range("K132").Value = 8705.454000000002
$ws.Range("L132").Value = 7543.956
$ws.Range("M132").Value = -6175.454000000002
$ws.Range("N132").Value = -12603.956
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
